# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" (copied from "2021-Q3" so it inherits the
#    same column layout / cell styles), positioned right before "总计",
#    then overwrite its header + two data rows with the 2022-Q1 fund data.
# 2. Prepend a new "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing rows down by one and renumbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the new "2022-Q1" detail sheet just before "总计".
# Worksheet variables captured before the Copy() below can end up
# pointing at the wrong (shifted) sheet position afterwards, so every
# sheet reference we need post-copy is re-looked-up by name.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q3")
$template.Copy($wb.Worksheets.Item("总计"))

$newSheet = $wb.Worksheets.Item("2021-Q3 (2)")
$newSheet.Name = "2022-Q1"
$newSheet = $wb.Worksheets.Item("2022-Q1")

# Header row
$newSheet.Cells.Item(1,2).Value = "基金代码"
$newSheet.Cells.Item(1,3).Value = "基金名称"
$newSheet.Cells.Item(1,4).Value = "基金规模"
$newSheet.Cells.Item(1,5).Value = "股票总仓位"
$newSheet.Cells.Item(1,6).Value = "仓位占比"
$newSheet.Cells.Item(1,7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1,8).Value = "仓位排名"

# Row 2 - fund 011815
$newSheet.Cells.Item(2,1).Value = 0
$newSheet.Cells.Item(2,2).Value = "011815"
$newSheet.Cells.Item(2,3).Value = "恒越优势精选混合型发起式证券投资基金"
$newSheet.Cells.Item(2,4).Value = "4.64"
$newSheet.Cells.Item(2,5).Value = "92.44"
$newSheet.Cells.Item(2,6).Value = "2.95"
$newSheet.Cells.Item(2,7).Value = "0.1369"
$newSheet.Cells.Item(2,8).Value = 10

# Row 3 - fund 008313
$newSheet.Cells.Item(3,1).Value = 1
$newSheet.Cells.Item(3,2).Value = "008313"
$newSheet.Cells.Item(3,3).Value = "光大保德信研究精选混合"
$newSheet.Cells.Item(3,4).Value = "2.64"
$newSheet.Cells.Item(3,5).Value = "88.98"
$newSheet.Cells.Item(3,6).Value = "4.02"
$newSheet.Cells.Item(3,7).Value = "0.1061"
$newSheet.Cells.Item(3,8).Value = 10

# ---------------------------------------------------------------------
# Step 2: prepend a "2022-Q1" row to the "总计" sheet (shift rows 2-4
# down to rows 3-5, keep the same per-row style by growing from the
# bottom up, then fill in the new top row). Re-fetch "总计" by name
# since the sheet Copy() above may have shifted worksheet positions.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")
$totals.Range("A4").Copy()
$totals.Range("A5").PasteSpecial(-4122)
$totals.Cells.Item(5,1).Value = 3
$totals.Cells.Item(5,2).Value = "2021-Q1"
$totals.Cells.Item(5,3).Value = 10
$totals.Cells.Item(5,4).Value = 1.23

$totals.Cells.Item(4,1).Value = 2
$totals.Cells.Item(4,2).Value = "2021-Q2"
$totals.Cells.Item(4,3).Value = 6
$totals.Cells.Item(4,4).Value = 0.39

$totals.Cells.Item(3,1).Value = 1
$totals.Cells.Item(3,2).Value = "2021-Q3"
$totals.Cells.Item(3,3).Value = 2
$totals.Cells.Item(3,4).Value = 0.24

$totals.Cells.Item(2,1).Value = 0
$totals.Cells.Item(2,2).Value = "2022-Q1"
$totals.Cells.Item(2,3).Value = 2
$totals.Cells.Item(2,4).Value = 0.24

# Restore the original active sheet / selection.
$wb.Worksheets.Item("2021-Q1").Activate()
